$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 1.8
$ws.Range("I3").Value = 3.8
$ws.Range("J3").Value = 1.02
$ws.Range("K3").Value = 21
$ws.Range("L3").Value = 1.11
$ws.Range("M3").Value = 6.5
$ws.Range("P3").Value = 1.22
$ws.Range("Q3").Value = 4
$ws.Range("R3").Value = 1.44
$ws.Range("S3").Value = 2.63
$ws.Range("U3").Value = 12
$ws.Range("AB3").Value = 12
$ws.Range("AE3").Value = 19
$ws.Range("J4").Value = 1.05
$ws.Range("K4").Value = 11
$ws.Range("K22").Value = 10
$ws.Range("N22").Value = 2.03
$ws.Range("O22").Value = 1.78
$ws.Range("G24").Value = 2.45
$ws.Range("H24").Value = 2.88
$ws.Range("I24").Value = 3.25
$ws.Range("J24").Value = 1.06
$ws.Range("K24").Value = 9.5
$ws.Range("R24").Value = 1.7
$ws.Range("S24").Value = 2.05
$ws.Range("U24").Value = 12
$ws.Range("W24").Value = 23
$ws.Range("X24").Value = 19
$ws.Range("Z24").Value = 9.5
$ws.Range("AB24").Value = 12
$ws.Range("AE24").Value = 11
$ws.Range("AF24").Value = 17
$ws.Range("AG24").Value = 12
$ws.Range("AH24").Value = 34
$ws.Range("G28").Value = 3.25
$ws.Range("I28").Value = 2.25
$ws.Range("N28").Value = 1.9
$ws.Range("O28").Value = 1.9
$ws.Range("R28").Value = 1.67
$ws.Range("S28").Value = 2.1
$ws.Range("AA28").Value = 6
$ws.Range("J38").Value = 1.03
$ws.Range("K38").Value = 17
